# PurchaseList.xlsx update
# Commit: "red and green led supplier info change due to availablity. run  BOM outjob"
#
# - LED (row 6) supplier part changed from OSRAM "LS L29K-G1J2-1-Z" to OSRAM "LS Q976-NR-1"
#   (description becomes a Super Red LED), with new stock/price.
# - LED (row 7) supplier part changed from OSRAM "LG L29K-F2J1-24-Z" to Kingbright "APT1608SGC"
#   (description becomes a Green water-clear LED), with new stock/price.
# - Several other rows' supplier stock quantities (column L) were refreshed from the
#   live BOM job.
# - Report created date/time (G24/H24) refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: LED -> Super Red (OSRAM) ---
$ws.Range("G6").Value = "Standard LEDs - SMD Super Red, 633nm 180mcd, 20mA"
$ws.Range("I6").Value = "LS Q976-NR-1"
$ws.Range("M6").Value = "720-LSQ976-NR-1"
$ws.Range("L6").Value = 378706
$ws.Range("R6").Value = 0.146

# --- Row 7: LED -> Green water clear (Kingbright) ---
$ws.Range("G7").Value = "Standard LEDs - SMD GREEN WATER CLEAR"
$ws.Range("H7").Value = "Kingbright"
$ws.Range("I7").Value = "APT1608SGC"
$ws.Range("M7").Value = "604-APT1608SGC"
$ws.Range("L7").Value = 131423
$ws.Range("R7").Value = 0.071

# --- Supplier stock (column L) refresh on other rows ---
$ws.Range("L2").Value = 310846
$ws.Range("L3").Value = 10285
$ws.Range("L4").Value = 84787
$ws.Range("L5").Value = 718242
$ws.Range("L8").Value = 190552
$ws.Range("L10").Value = 18014
$ws.Range("L14").Value = 987320
$ws.Range("L17").Value = 25904
$ws.Range("L18").Value = 216662
$ws.Range("L21").Value = 23179
$ws.Range("L22").Value = 2369

# --- Report created date/time ---
$ws.Range("G24").Value = "31/10/2017"
$ws.Range("H24").Value = "09:26:15"

# --- Keep the hyperlink display text in sync with the new cell text ---
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$H$7') {
        $hl.TextToDisplay = "'Kingbright"
    } elseif ($addr -eq '$I$6') {
        $hl.TextToDisplay = "'LS Q976-NR-1"
    } elseif ($addr -eq '$I$7') {
        $hl.TextToDisplay = "'APT1608SGC"
    } elseif ($addr -eq '$M$6') {
        $hl.TextToDisplay = "'720-LSQ976-NR-1"
    } elseif ($addr -eq '$M$7') {
        $hl.TextToDisplay = "'604-APT1608SGC"
    }
}
